$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.386.77"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.067.20"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "234.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("E7").Value = "  -0.02%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "57.51"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.05%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.398"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.95%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0773"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "2.370.87"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  -0.27%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.75"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.776"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.19"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "2.067.46"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").Value = "37.337.64"
$ws.Range("E18").Value = "  -0.69%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.26"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.49%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "69.65"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  +0.32%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "226.22"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("E23").Value = "  -0.05%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").Value = "  -2.01%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "167.01"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.87%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.86"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.22%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.10"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("E30").Value = "  +0.47%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.117"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.54"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0618"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.55"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.39%  "

$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  -3.41%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -4.61%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.67%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0965"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.96%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "98.01"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.84%  "

$ws.Range("D43").Value = "1.483.95"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("E44").Value = "  +0.93%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.17"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.06"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -11.76%  "

$ws.Range("E47").Value = "  +0.07%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "15.32"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").Value = "2.257.66"
$ws.Range("E51").Value = "  +0.12%  "
